# CV-Michel_de_Bree.EN.docx edits
# 1) "With over 20+ years ... I am" / "equiped to help you ..." -> "With 15+ years ... I am equiped" / "to help you ..."
# 2) "JUnit / Mockito / Cucumber / Gherkin / Wiremock" -> "JUnit / Mockito / Wiremock"
# 3) "Java 8, Spring Boot, ReactJS, Docker, Drools" -> "Java 8, Spring Boot, Hibernate, ReactJS, Docker, Drools"
# 4) Keywords line: tone down the keyword list / rename some entries

$d = $word.ActiveDocument

function Split-RunAt($doc, $pos) {
    # Forces a run boundary at an absolute character offset by adding and
    # immediately removing a zero-length bookmark there. This leaves no
    # residual formatting marks (unlike toggling a character property).
    $b = $doc.Range($pos, $pos)
    $name = "tmp_split_" + $pos
    $doc.Bookmarks.Add($name, $b) | Out-Null
    $doc.Bookmarks($name).Delete()
}

# --- Edit 1: intro paragraph -------------------------------------------------
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.StartsWith("With over 20+ years")) {
        $pStart = $p.Range.Start
        $oldFull = "With over 20+ years of experience in complex and dynamic IT landscapes, I am equiped to help you build your automated solutions in an agile and sustainable way."
        $newFull = "With 15+ years of experience in complex and dynamic IT landscapes, I am equiped to help you build your automated solutions in an agile and sustainable way."
        $newRun1 = "With 15+ years of experience in complex and dynamic IT landscapes, I am equiped"

        $fullEnd = $pStart + $oldFull.Length
        $fullRange = $d.Range($pStart, $fullEnd)
        $fullRange.Text = $newFull

        $split1 = $pStart + $newRun1.Length
        $split2 = $split1 + 1
        Split-RunAt $d $split1
        Split-RunAt $d $split2
        break
    }
}

# --- Edit 2: JUnit / Mockito / Cucumber / Gherkin / Wiremock ---------------
$d.Content.Find.Execute("JUnit / Mockito / Cucumber / Gherkin / Wiremock", $true, $false, $false, $false, $false, `
    $true, 1, $false, "JUnit / Mockito / Wiremock", 2) | Out-Null

# --- Edit 3: Java 8, Spring Boot, ReactJS, Docker, Drools -------------------
$d.Content.Find.Execute("Java 8, Spring Boot, ReactJS, Docker, Drools", $true, $false, $false, $false, $false, `
    $true, 1, $false, "Java 8, Spring Boot, Hibernate, ReactJS, Docker, Drools", 2) | Out-Null

# --- Edit 4: Keywords list on the Rijkswaterstaat / rule engine project ----
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "Keywords:*Full Stack Developer, Agile, Scrum, DevOps, SAFE*") {
        $pStart = $p.Range.Start
        $prefix = "Keywords:"
        $restStart = $pStart + $prefix.Length
        $restOldText = " Full Stack Developer, Agile, Scrum, DevOps, SAFE, Java, ReactJS, REST, JSON, OAuth, Swagger/OpenAPI, Docker, Drools, Postgres, Continuous Integration, JIRA, Confluence, Bitbucket, Maven, Jenkins"
        $restNewText = " Full Stack Developer, Scrum, SAFE, Java, Spring Boot, Hibernate, ReactJS, JUnit, Mockito, REST, JSON, OAuth, Swagger/OpenAPI, Docker, Drools, Postgres, CI/CD, JIRA, Confluence, Bitbucket, Maven, Jenkins"
        $restEnd = $restStart + $restOldText.Length

        $r = $d.Range($restStart, $restEnd)
        $r.Text = $restNewText

        # Re-establish the original run boundaries (space / phrase / space / phrase / space / phrase)
        $splitOffsets = @(1, 90, 91, 159, 160)
        foreach ($off in $splitOffsets) {
            $splitPos = $restStart + $off
            Split-RunAt $d $splitPos
        }
        break
    }
}
